$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 282, pushing existing data (rows 282..391)
# down to rows 283..392. Excel's native row-insert also fixes up the
# sheet's used range / dimension automatically.
$ws.Rows.Item(282).Insert()

# Populate the newly-inserted row 282 with the new record. Most columns
# mirror the row that used to occupy 282 (now at 283); only D, J, L, M, P
# carry new values for this record.
$ws.Cells.Item(282, 1).Value = 3
$ws.Cells.Item(282, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(282, 3).Value = "Coquimbo"
$ws.Cells.Item(282, 4).Value = 44795
$ws.Cells.Item(282, 5).Value = 5
$ws.Cells.Item(282, 6).Value = 100112012
$ws.Cells.Item(282, 7).Value = "Espinaca"
$ws.Cells.Item(282, 8).Value = "Sin especificar"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 280
$ws.Cells.Item(282, 11).Value = 4000
$ws.Cells.Item(282, 12).Value = 4200
$ws.Cells.Item(282, 13).Value = 4093
$ws.Cells.Item(282, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(282, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(282, 16).Value = 1364
$ws.Cells.Item(282, 17).Value = 3
$ws.Cells.Item(282, 18).Value = "Hortaliza"
